# Generate Report for Handoff
# Update the "Latest Handoff Date/Datetime" values for the files that were
# re-handed-off, across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Rows 7, 10-16 on every sheet correspond to the same set of files
# (927eb88a..., 17a40e05..., 27999d32..., 6bb99208..., b0972425...,
#  da380b96..., e6bb96d2..., f953f76f...) whose handoff timestamp was
# refreshed by this report generation run.
$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

foreach ($r in $rows) {
    $overview.Range("D$r").Value = "2016-03-22 04:29:14"
    $zhcn.Range("E$r").Value     = "2016-03-22 04:29:10"
    $dede.Range("E$r").Value     = "2016-03-22 04:29:14"
}
